# Apply forecast data updates (Optuna attempt - revert to original values)

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---

# Row 2 (W8)
$wsForecast.Range("D2").Value = 375
$wsForecast.Range("H2").Value = 5.11

# Row 3 (W9)
$wsForecast.Range("D3").Value = 393
$wsForecast.Range("H3").Value = 3.92
$wsForecast.Range("L3").Value = 0.8

# Row 4 (W10)
$wsForecast.Range("D4").Value = 364
$wsForecast.Range("H4").Value = 3.16
$wsForecast.Range("L4").Value = 0.85

# Row 5 (W11)
$wsForecast.Range("D5").Value = 356
$wsForecast.Range("H5").Value = 2.2
$wsForecast.Range("L5").Value = 1.06

# Row 6 (W12)
$wsForecast.Range("H6").Value = 1.17
$wsForecast.Range("L6").Value = 1.03

# Row 7 (W13)
$wsForecast.Range("H7").Value = 0.19
$wsForecast.Range("I7").Value = "High"
$wsForecast.Range("L7").Value = 1.04

# Row 8 (W14)
$wsForecast.Range("L8").Value = 1.04

# Row 9 (W15)
$wsForecast.Range("L9").Value = 0.89

# Row 10 (W16)
$wsForecast.Range("L10").Value = 1.16

# Row 11 (W17)
$wsForecast.Range("L11").Value = 1

# Row 12 (W18)
$wsForecast.Range("L12").Value = 0.9

# Row 13 (W19)
$wsForecast.Range("L13").Value = 1.07

# Row 14 (W20)
$wsForecast.Range("L14").Value = 1.11

# Row 15 (W21)
$wsForecast.Range("L15").Value = 1.05

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.14

# --- Summary sheet ---
# These cells store numeric-looking values as text (as in the source file),
# so prefix with an apostrophe to keep Excel from re-typing them as numbers.
$wsSummary.Range("B9").Value = "'5318"
$wsSummary.Range("B10").Value = "'2858"
$wsSummary.Range("B11").Value = "'1488"
$wsSummary.Range("B12").Value = "'393"
